# Update countries & provincias Spain
# Applies the 31-Aug-2020 20:00 -> 21:17 data refresh to the "Pais" sheet:
#  - bumps the "last updated" timestamp string
#  - updates case counters for the countries whose figures changed
#  - re-sorts the rows whose case totals crossed a neighbour (which also
#    shuffles which country name now belongs on those rows)

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# --- Last updated timestamp (row 1) ---------------------------------------
$ws.Range("A1").Value = "Datos actualizados a 31 de Agosto de 2020 a las 21:17"

# --- Simple in-place numeric refreshes (country/rank unchanged) ----------
# row 4: Estados Unidos
$ws.Range("B4").Value = 6194854
$ws.Range("C4").Value = 21618
$ws.Range("D4").Value = 3433746
$ws.Range("E4").Value = 2573601
$ws.Range("G4").Value = 283
$ws.Range("H4").Value = 187507

# row 5: Brasil
$ws.Range("B5").Value = 3866157
$ws.Range("C5").Value = 3846
$ws.Range("E5").Value = 713622
$ws.Range("G5").Value = 80
$ws.Range("H5").Value = 120976

# row 20: Francia
$ws.Range("D20").Value = 86469
$ws.Range("E20").Value = 163921

# row 23: Alemania
$ws.Range("B23").Value = 244703
$ws.Range("C23").Value = 1408
$ws.Range("E23").Value = 17490
$ws.Range("G23").Value = 7
$ws.Range("H23").Value = 9371

# row 24: Irak
$ws.Range("B24").Value = 234934
$ws.Range("C24").Value = 3757
$ws.Range("D24").Value = 176602
$ws.Range("E24").Value = 51290
$ws.Range("G24").Value = 83
$ws.Range("H24").Value = 7042

# row 97: Guayana Francesa
$ws.Range("B97").Value = 9115
$ws.Range("C97").Value = 39
$ws.Range("D97").Value = 8654
$ws.Range("E97").Value = 402
$ws.Range("G97").Value = 1
$ws.Range("H97").Value = 59

# row 127: Sri Lanka
$ws.Range("B127").Value = 3049
$ws.Range("C127").Value = 37
$ws.Range("E127").Value = 169

# row 133: Siria
$ws.Range("B133").Value = 2765
$ws.Range("C133").Value = 62
$ws.Range("D133").Value = 629
$ws.Range("E133").Value = 2024
$ws.Range("G133").Value = 3
$ws.Range("H133").Value = 112

# row 134: Angola
$ws.Range("B134").Value = 2654
$ws.Range("C134").Value = 30
$ws.Range("D134").Value = 1071
$ws.Range("E134").Value = 1475
$ws.Range("G134").Value = 1
$ws.Range("H134").Value = 108

# row 135: Sudan del Sur
$ws.Range("B135").Value = 2527
$ws.Range("C135").Value = 8
$ws.Range("E135").Value = 1190

# row 144: Aruba
$ws.Range("B144").Value = 2006
$ws.Range("C144").Value = 9
$ws.Range("D144").Value = 768
$ws.Range("E144").Value = 1228

# row 145: Yemen
$ws.Range("B145").Value = 1958
$ws.Range("C145").Value = 5
$ws.Range("D145").Value = 1131
$ws.Range("E145").Value = 261
$ws.Range("G145").Value = 2
$ws.Range("H145").Value = 566

# row 150: Reunion
$ws.Range("B150").Value = 1679
$ws.Range("C150").Value = 45
$ws.Range("E150").Value = 790

# row 152: Republica de Chipre
$ws.Range("B152").Value = 1488
$ws.Range("C152").Value = 1
$ws.Range("E152").Value = 329

# row 182: Eritrea
$ws.Range("B182").Value = 319
$ws.Range("C182").Value = 1
$ws.Range("E182").Value = 35

# --- Re-sorted block (rows 116-119 & 121-122) -----------------------------
# Cuba's update (3973 -> 4032 total cases) pushes it above Ruanda, Surinam
# and Congo; Mozambique's update (3821 -> 3916) pushes it above Cabo Verde.
# Rewrite the country name + full data row for each affected position.

# row 116 becomes Cuba (was Ruanda)
$ws.Range("A116").Value = "Cuba"
$ws.Range("B116").Value = 4032
$ws.Range("C116").Value = 59
$ws.Range("D116").Value = 3378
$ws.Range("E116").Value = 560
$ws.Range("F116").Value = 0
$ws.Range("G116").Value = 0
$ws.Range("H116").Value = 94

# row 117 becomes Ruanda (was Surinam) - unchanged figures, just moved down
$ws.Range("A117").Value = "Ruanda"
$ws.Range("B117").Value = 4020
$ws.Range("C117").Value = 0
$ws.Range("D117").Value = 1918
$ws.Range("E117").Value = 2086
$ws.Range("F117").Value = 0
$ws.Range("G117").Value = 0
$ws.Range("H117").Value = 16

# row 118 becomes Surinam (was Congo) - unchanged figures, just moved down
$ws.Range("A118").Value = "Surinam"
$ws.Range("B118").Value = 4009
$ws.Range("C118").Value = 0
$ws.Range("D118").Value = 3073
$ws.Range("E118").Value = 869
$ws.Range("F118").Value = 0
$ws.Range("G118").Value = 0
$ws.Range("H118").Value = 67

# row 119 becomes Congo (was Cuba) - unchanged figures, just moved down
$ws.Range("A119").Value = "Congo"
$ws.Range("B119").Value = 3979
$ws.Range("C119").Value = 0
$ws.Range("D119").Value = 1742
$ws.Range("E119").Value = 2159
$ws.Range("F119").Value = 0
$ws.Range("G119").Value = 0
$ws.Range("H119").Value = 78

# row 120 (Eslovaquia) unchanged

# row 121 becomes Mozambique (was Cabo Verde)
$ws.Range("A121").Value = "Mozambique"
$ws.Range("B121").Value = 3916
$ws.Range("C121").Value = 95
$ws.Range("D121").Value = 2170
$ws.Range("E121").Value = 1723
$ws.Range("F121").Value = 0
$ws.Range("G121").Value = 0
$ws.Range("H121").Value = 23

# row 122 becomes Cabo Verde (was Mozambique)
$ws.Range("A122").Value = "Cabo Verde"
$ws.Range("B122").Value = 3884
$ws.Range("C122").Value = 32
$ws.Range("D122").Value = 2916
$ws.Range("E122").Value = 928
$ws.Range("F122").Value = 0
$ws.Range("G122").Value = 0
$ws.Range("H122").Value = 40
